# Auto-generated Excel COM-interop script to update cryptos.xlsx data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.455.98"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "2.221.61"
$ws.Range("E3").Value = "  -4.42%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'297.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.41%  "
$ws.Range("D6").Value = "'82.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.85%  "
$ws.Range("E7").Value = "  -3.22%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.471"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.49%  "
$ws.Range("E10").Value = "  -4.26%  "
$ws.Range("D11").Value = "'29.79"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("D12").Value = "'46.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -11.14%  "
$ws.Range("E13").Value = "  -2.60%  "
$ws.Range("D15").Value = "'6.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.63%  "
$ws.Range("D16").Value = "'14.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.71%  "
$ws.Range("D17").Value = "2.220.02"
$ws.Range("E17").Value = "  -4.62%  "
$ws.Range("D18").Value = "'0.717"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.59%  "
$ws.Range("D19").Value = "39.379.49"
$ws.Range("E19").Value = "  -1.21%  "
$ws.Range("D20").Value = "0.0₃0876"
$ws.Range("E20").Value = "  -2.77%  "
$ws.Range("D21").Value = "'5.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.01%  "
$ws.Range("D22").Value = "'64.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.72%  "
$ws.Range("D23").Value = "'10.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.19%  "
$ws.Range("D24").Value = "'230.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.76%  "
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("E26").Value = "  -5.47%  "
$ws.Range("D27").Value = "'1.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("D28").Value = "'22.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.98%  "
$ws.Range("E29").Value = "  +1.89%  "
$ws.Range("D30").Value = "'9.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "'32.04"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.32%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "'149.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.51%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("D34").Value = "'4.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.85%  "
$ws.Range("D35").Value = "'0.0699"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.70%  "
$ws.Range("D36").Value = "'2.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.01%  "
$ws.Range("D37").Value = "'0.111"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.18%  "
$ws.Range("D38").Value = "'15.73"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.40%  "
$ws.Range("D39").Value = "'0.0963"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.14%  "
$ws.Range("E40").Value = "  -3.15%  "
$ws.Range("D41").Value = "'1.67"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.08%  "
$ws.Range("D42").Value = "'3.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.32%  "
$ws.Range("D43").Value = "1.914.34"
$ws.Range("E43").Value = "  -1.41%  "
$ws.Range("D45").Value = "'0.0261"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.67%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'16.49"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.83%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'9.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.47%  "
$ws.Range("D48").Value = "'2.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.89%  "
$ws.Range("D49").Value = "2.426.53"
$ws.Range("E49").Value = "  -4.98%  "
$ws.Range("D50").Value = "'71.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.72%  "
$ws.Range("D51").Value = "'88.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.72%  "
